# Applies the "wrapping up test file audit" edit:
#  - Remove the stray "Sheet" row (row 16: "Sheet", 3, 4) from the
#    optimization_parameters sheet, shifting the simulation_timepoints
#    row up from 17 to 16.
#  - Move the active/selected tab from optimization_parameters to
#    optimization_diagnostics (the last sheet).

$wb = $excel.ActiveWorkbook

$paramsSheet = $wb.Worksheets.Item("optimization_parameters")
$paramsSheet.Rows.Item(16).Delete()

$diagSheet = $wb.Worksheets.Item("optimization_diagnostics")
$diagSheet.Activate()
$diagSheet.Select()
